$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing it to remain a text string
# (so numeric-looking strings like "1.010" do not get coerced into numbers),
# then restore the cell to its original (default/"Normal") style so no
# stray formatting is introduced.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "27.445.98"
$ws.Range("E2").Value = "  +1.79%  "

# Row 3
$ws.Range("D3").Value = "1.858.36"
$ws.Range("E3").Value = "  +0.82%  "

# Row 4
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
Set-TextValue $ws.Range("D5") "315.54"
$ws.Range("E5").Value = "  +2.25%  "

# Row 6
Set-TextValue $ws.Range("D6") "1.010"
$ws.Range("E6").Value = "  -0.08%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.4768"
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.3796"
$ws.Range("E8").Value = "  +3.25%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.07304"
$ws.Range("E9").Value = "  +1.40%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.9295"

# Row 12
Set-TextValue $ws.Range("D12") "0.07787"
$ws.Range("E12").Value = "  +0.83%  "

# Row 13
$ws.Range("D13").Value = "1.859.43"
$ws.Range("E13").Value = "  +1.27%  "

# Row 14
Set-TextValue $ws.Range("D14") "5.439"
$ws.Range("E14").Value = "  +0.91%  "

# Row 15
Set-TextValue $ws.Range("D15") "6.543"
$ws.Range("E15").Value = "  +1.63%  "

# Row 16
Set-TextValue $ws.Range("D16") "90.18"
$ws.Range("E16").Value = "  +1.54%  "

# Row 17
Set-TextValue $ws.Range("D17") "1.012"

# Row 18
Set-TextValue $ws.Range("D18") "0.000008812"
$ws.Range("E18").Value = "  +2.05%  "

# Row 19
Set-TextValue $ws.Range("D19") "1.009"
$ws.Range("E19").Value = "  -0.28%  "

# Row 20
$ws.Range("D20").Value = "27.519.26"
$ws.Range("E20").Value = "  +2.10%  "

# Row 21
$ws.Range("E21").Value = "  +0.53%  "

# Row 22
Set-TextValue $ws.Range("D22") "5.093"
$ws.Range("E22").Value = "  +0.55%  "

# Row 23
Set-TextValue $ws.Range("D23") "10.68"
$ws.Range("E23").Value = "  +0.53%  "

# Row 24
Set-TextValue $ws.Range("D24") "1.944"
$ws.Range("E24").Value = "  -0.01%  "

# Row 25
Set-TextValue $ws.Range("D25") "154.81"
$ws.Range("E25").Value = "  +1.56%  "

# Row 26
Set-TextValue $ws.Range("D26") "18.44"
$ws.Range("E26").Value = "  +1.34%  "

# Row 27
Set-TextValue $ws.Range("D27") "1.998"
$ws.Range("E27").Value = "  -0.58%  "

# Row 28
Set-TextValue $ws.Range("D28") "115.32"
$ws.Range("E28").Value = "  +0.95%  "

# Row 29
Set-TextValue $ws.Range("D29") "4.936"
$ws.Range("E29").Value = "  -0.43%  "

# Row 30
Set-TextValue $ws.Range("D30") "0.08889"
$ws.Range("E30").Value = "  +0.38%  "

# Row 31
Set-TextValue $ws.Range("D31") "3.331"
$ws.Range("E31").Value = "  +1.14%  "

# Row 32
Set-TextValue $ws.Range("D32") "1.202"
$ws.Range("E32").Value = "  +2.44%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.7527"
$ws.Range("E33").Value = "  +1.85%  "

# Row 34
Set-TextValue $ws.Range("D34") "4.576"
$ws.Range("E34").Value = "  +1.78%  "

# Row 35
Set-TextValue $ws.Range("D35") "2.693"
$ws.Range("E35").Value = "  +0.17%  "

# Row 36
Set-TextValue $ws.Range("D36") "1.125"
$ws.Range("E36").Value = "  +1.44%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.02042"
$ws.Range("E37").Value = "  +4.32%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.5537"
$ws.Range("E38").Value = "  +5.61%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.05274"
$ws.Range("E39").Value = "  +0.41%  "

# Row 40
Set-TextValue $ws.Range("D40") "2.987"
$ws.Range("E40").Value = "  +0.84%  "

# Row 41
Set-TextValue $ws.Range("D41") "7.014"
$ws.Range("E41").Value = "  +0.27%  "

# Row 42
$ws.Range("E42").Value = "  +3.32%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.1513"
$ws.Range("E43").Value = "  +0.25%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.4860"
$ws.Range("E44").Value = "  +2.75%  "

# Row 45
Set-TextValue $ws.Range("D45") "10.62"
$ws.Range("E45").Value = "  +0.24%  "

# Row 46
Set-TextValue $ws.Range("D46") "1.011"
$ws.Range("E46").Value = "  -0.12%  "

# Row 47
$ws.Range("E47").Value = "  +3.72%  "

# Row 48
Set-TextValue $ws.Range("D48") "103.05"
$ws.Range("E48").Value = "  +1.33%  "

# Row 49
Set-TextValue $ws.Range("D49") "67.37"
$ws.Range("E49").Value = "  +2.71%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.06098"
$ws.Range("E50").Value = "  +0.52%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.9133"
$ws.Range("E51").Value = "  +2.86%  "
